$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 (header row): F1 was a boolean (TRUE) placeholder, now becomes the
# text header for a new "Blocked" column: "Blocked Mood:"
$ws.Range("F1").Value = "Blocked Mood:"

# --- Row 4: Amr Elsayed Elhenawy's account record is replaced by a brand new
# account record for Ahmed Nasr.
$ws.Range("A4").Value = "ahmed.nasr"
$ws.Range("B4").Value = "Ahmad Nasr"
$ws.Range("C4").Value = "01125697852"
$ws.Range("D4").Value = "ahmed.nasr@gmail.com"
$ws.Range("E4").Value = "ZsnvwMzhi123#"
$ws.Range("F4").Value = $false

# --- Row 5 (new row): Amr Elhenawy's account is re-added afterwards with an
# updated username/phone/email/password.
$ws.Range("A5").Value = "amr.elhenawy"
$ws.Range("B5").Value = "Amr Elhenawy"
$ws.Range("C5").Value = "01234567891"
$ws.Range("D5").Value = "amr.elhenawy123@gmail.com"
$ws.Range("E5").Value = "zni.vosvmzdb123@tnzro.xln"
$ws.Range("F5").Value = $false

# The new row was not populated from the worksheet's per-column formatting
# (rows 2-4 carry s="1"/s="2" styles inherited from the column defaults), so
# normalize row 5 back to the workbook's default "Normal" style.
$ws.Range("A5:F5").Style = "Normal"
